$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "트리거가 필요한 실험에서 결과의 신뢰성을 높이는 방법"
$ws.Range("E3").Value = "https://lumiamitie.github.io/data/counterfactual-logging/"

$ws.Range("D9").Value = "한국 대학 vs. 해외 대학 (2)"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/korean-uni-vs-siai-2/#utm_source=rss&utm_medium=rss&utm_campaign=korean-uni-vs-siai-2"

$ws.Range("D44").Value = "Intel의 Neuromorphic Chip - Loihi 2"
$ws.Range("E44").Value = "https://engineering-ladder.tistory.com/98"
